# Trade #38 closed at 2026-02-17 08:33:06 - unknown UNKNOWN +0.000%
#
# This script updates the workbook to reflect the closing of Trade #38
# (MarketMaking strategy) and the resulting roll-up changes to the
# Summary and Strategy Status sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet - overall portfolio roll-up numbers
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.66    # Current Capital
$summary.Range("B4").Value = -0.34      # Total P&L $
$summary.Range("B5").Value = -0.18      # Total P&L %
$summary.Range("B6").Value = 38         # Total Trades
$summary.Range("B7").Value = 13         # Winning Trades
$summary.Range("B9").Value = 34.21      # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking strategy row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.66       # Capital
$status.Range("D4").Value = 38          # Trades
$status.Range("E4").Value = -0.34       # P&L $
$status.Range("F4").Value = -0.34       # P&L %
$status.Range("G4").Value = 34.21       # Win Rate %

# ---------------------------------------------------------------
# All Trades & MarketMaking sheets - Trade #38 (row 39) now closed
# ---------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G39").Value = 0.74                    # Exit Price
    $ws.Range("H39").Value = "CLOSED"                # Status
    $ws.Range("I39").Value = 10.4478                 # P&L %
    $ws.Range("J39").Value = 0.07000000000000001     # P&L $
    $ws.Range("K39").Value = 99.66                   # Capital After
    $ws.Range("P39").Value = "early_exit"             # Exit Reason
    $ws.Range("Q39").Value = 0.13                    # Duration (min)
}
